$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Style J21 (Komi text) in bold green - matches new font/xf added to styles.xml
$ws.Range("J21").Font.Bold = $true
$ws.Range("J21").Font.Color = 5287936

# Style J11 (Determining Winner text) in bold orange - matches new font/xf added to styles.xml
$ws.Range("J11").Font.Bold = $true
$ws.Range("J11").Font.Color = 49407

# Add new rows 31-33 with "B" labels extracted into the rules area (G:I columns)
$ws.Range("G31:I33").Value = "B"

# Move the active selection to J11
$ws.Range("J11").Select()
